# COVID19HNHistory.xlsx - "Add files via upload"
# A new daily snapshot column (BA, date serial 43973 = the day after the
# existing AZ column's 43972) is appended to the "Evol. No. casos" sheet,
# with one case count per department plus a total (SUM) row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: next day's date serial, right after AZ1 (43972)
$ws.Range("BA1").Value = 43973

# Per-department case counts for the new date (rows 2-19)
$counts = @{
    2  = 109
    3  = 80
    4  = 42
    5  = 31
    6  = 5
    7  = 2184
    8  = 23
    9  = 705
    10 = 0
    11 = 12
    12 = 0
    13 = 8
    14 = 12
    15 = 1
    16 = 9
    17 = 26
    18 = 45
    19 = 185
}

foreach ($row in $counts.Keys) {
    $ws.Cells.Item($row, 53).Value = $counts[$row]
}

# Totals row mirrors the AZ20 SUM formula for the new column
$ws.Range("BA20").Formula = "=SUM(BA2:BA19)"

# Leave the selection on the new column, matching where the author's
# cursor ended up (BA2:BA20, landed on the totals cell)
$ws.Range("BA2:BA20").Select() | Out-Null
